$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.77516109082218
$ws.Range("C2").Value = 12.00691989816932
$ws.Range("E2").Value = 15.91959507680913
$ws.Range("F2").Value = 38.67410743912644
$ws.Range("G2").Value = 3.65945006192706
$ws.Range("J2").Value = 8.232652121800822
$ws.Range("L2").Value = 12.31552533303935
$ws.Range("M2").Value = 16.96736302495654
$ws.Range("O2").Value = 25.17923915579559
$ws.Range("B3").Value = 16.28589118250084
$ws.Range("C3").Value = 11.87400362820946
$ws.Range("E3").Value = 15.9596418490704
$ws.Range("F3").Value = 38.77913314024681
$ws.Range("G3").Value = 3.661644842525671
$ws.Range("J3").Value = 8.233091700986588
$ws.Range("L3").Value = 12.30519580041614
$ws.Range("M3").Value = 16.84823663255789
$ws.Range("O3").Value = 25.294088368547
$ws.Range("B4").Value = 15.97954520309201
$ws.Range("C4").Value = 11.79145150709885
$ws.Range("E4").Value = 15.98588193838138
$ws.Range("F4").Value = 38.85345887512945
$ws.Range("G4").Value = 3.663063612318203
$ws.Range("J4").Value = 8.233665928874686
$ws.Range("L4").Value = 12.30010767453179
$ws.Range("M4").Value = 16.77654884899353
$ws.Range("O4").Value = 25.37114008423291
$ws.Range("B5").Value = 15.85339536486431
$ws.Range("C5").Value = 11.75759702162447
$ws.Range("E5").Value = 15.99699094863865
$ws.Range("F5").Value = 38.88621447878486
$ws.Range("G5").Value = 3.663659725267994
$ws.Range("J5").Value = 8.233976828167259
$ws.Range("L5").Value = 12.2983514561307
$ws.Range("M5").Value = 16.74772350628098
$ws.Range("O5").Value = 25.40417718201254
$ws.Range("B6").Value = 15.83237459157837
$ws.Range("C6").Value = 11.75196326592072
$ws.Range("E6").Value = 15.99886073646919
$ws.Range("F6").Value = 38.8918023243335
$ws.Range("G6").Value = 3.66375979536375
$ws.Range("J6").Value = 8.234033107981254
$ws.Range("L6").Value = 12.29807904925634
$ws.Range("M6").Value = 16.74296114054424
$ws.Range("O6").Value = 25.40976174233442
$ws.Range("B7").Value = 15.97784896446716
$ws.Range("C7").Value = 11.79099576916534
$ws.Range("E7").Value = 15.98603007316471
$ws.Range("F7").Value = 38.8538906484345
$ws.Range("G7").Value = 3.663071578939126
$ws.Range("J7").Value = 8.233669809935378
$ws.Range("L7").Value = 12.30008270273119
$ws.Range("M7").Value = 16.77615850093975
$ws.Range("O7").Value = 25.37157900909268
$ws.Range("B8").Value = 16.60779637786373
$ws.Range("C8").Value = 11.96129926553821
$ws.Range("E8").Value = 15.9330609473871
$ws.Range("F8").Value = 38.70827407552061
$ws.Range("G8").Value = 3.660192086566
$ws.Range("J8").Value = 8.232740748790194
$ws.Range("L8").Value = 12.31170452592896
$ws.Range("M8").Value = 16.92599924599533
$ws.Range("O8").Value = 25.21748002454568
$ws.Range("B9").Value = 17.78877577309483
$ws.Range("C9").Value = 12.28681307912102
$ws.Range("E9").Value = 15.84225641330729
$ws.Range("F9").Value = 38.50108241747566
$ws.Range("G9").Value = 3.655107448980128
$ws.Range("J9").Value = 8.233317596673572
$ws.Range("L9").Value = 12.34436282204513
$ws.Range("M9").Value = 17.23041000597937
$ws.Range("O9").Value = 24.96736074140238
$ws.Range("B10").Value = 18.61449249295694
$ws.Range("C10").Value = 12.51952683236968
$ws.Range("E10").Value = 15.78346200830314
$ws.Range("F10").Value = 38.39699340490941
$ws.Range("G10").Value = 3.651710715607415
$ws.Range("J10").Value = 8.23518288737851
$ws.Range("L10").Value = 12.37425560757147
$ws.Range("M10").Value = 17.4591640086743
$ws.Range("O10").Value = 24.81564525509172
$ws.Range("B11").Value = 18.97942095693964
$ws.Range("C11").Value = 12.62373126719485
$ws.Range("E11").Value = 15.75842473541385
$ws.Range("F11").Value = 38.36015755606545
$ws.Range("G11").Value = 3.650238269201971
$ws.Range("J11").Value = 8.23634001576136
$ws.Range("L11").Value = 12.38910830852441
$ws.Range("M11").Value = 17.56404597948362
$ws.Range("O11").Value = 24.75364888375823
$ws.Range("B12").Value = 19.11595549413624
$ws.Range("C12").Value = 12.66293081199452
$ws.Range("E12").Value = 15.7491887231954
$ws.Range("F12").Value = 38.34772497108819
$ws.Range("G12").Value = 3.649691093184939
$ws.Range("J12").Value = 8.236822177126783
$ws.Range("L12").Value = 12.3949104810612
$ws.Range("M12").Value = 17.60385512140327
$ws.Range("O12").Value = 24.73118714366218
$ws.Range("B13").Value = 19.08662574081998
$ws.Range("C13").Value = 12.6545004133027
$ws.Range("E13").Value = 15.75116697462128
$ws.Range("F13").Value = 38.35033504923301
$ws.Range("G13").Value = 3.649808475296203
$ws.Range("J13").Value = 8.236716385245956
$ws.Range("L13").Value = 12.39365301520652
$ws.Range("M13").Value = 17.59527779591802
$ws.Range("O13").Value = 24.73597945483889
$ws.Range("B14").Value = 18.99068755712034
$ws.Range("C14").Value = 12.62696157987166
$ws.Range("E14").Value = 15.75765997537979
$ws.Range("F14").Value = 38.35910430760885
$ws.Range("G14").Value = 3.650193044423055
$ws.Range("J14").Value = 8.236378804271292
$ws.Range("L14").Value = 12.38958210909011
$ws.Range("M14").Value = 17.5673193909466
$ws.Range("O14").Value = 24.75178057075557
$ws.Range("B15").Value = 18.93170369535233
$ws.Range("C15").Value = 12.61005869508706
$ws.Range("E15").Value = 15.76166902111039
$ws.Range("F15").Value = 38.36467332167224
$ws.Range("G15").Value = 3.650429957908999
$ws.Range("J15").Value = 8.236177742178304
$ws.Range("L15").Value = 12.38711163447021
$ws.Range("M15").Value = 17.55020536771143
$ws.Range("O15").Value = 24.76159154754178
$ws.Range("B16").Value = 18.59041798880476
$ws.Range("C16").Value = 12.5126815676366
$ws.Range("E16").Value = 15.78513257991364
$ws.Range("F16").Value = 38.3996126859077
$ws.Range("G16").Value = 3.651808402701747
$ws.Range("J16").Value = 8.235113442483035
$ws.Range("L16").Value = 12.37330996975887
$ws.Range("M16").Value = 17.45232406925012
$ws.Range("O16").Value = 24.8198385921479
$ws.Range("B17").Value = 18.37822328843305
$ws.Range("C17").Value = 12.45250331425422
$ws.Range("E17").Value = 15.7999638722937
$ws.Range("F17").Value = 38.42374325642776
$ws.Range("G17").Value = 3.652672627667356
$ws.Range("J17").Value = 8.234539289072611
$ws.Range("L17").Value = 12.36516254900774
$ws.Range("M17").Value = 17.3924692021408
$ws.Range("O17").Value = 24.85737314377937
$ws.Range("B18").Value = 18.25517418913766
$ws.Range("C18").Value = 12.41773556973741
$ws.Range("E18").Value = 15.80865529813872
$ws.Range("F18").Value = 38.43861200739742
$ws.Range("G18").Value = 3.653176557229467
$ws.Range("J18").Value = 8.234238121679516
$ws.Range("L18").Value = 12.36059459333869
$ws.Range("M18").Value = 17.35812103291461
$ws.Range("O18").Value = 24.87962243451912
$ws.Range("B19").Value = 18.2133438715719
$ws.Range("C19").Value = 12.40593790070581
$ws.Range("E19").Value = 15.81162571267616
$ws.Range("F19").Value = 38.44381611189615
$ws.Range("G19").Value = 3.653348357278104
$ws.Range("J19").Value = 8.234141156421067
$ws.Range("L19").Value = 12.35906834447557
$ws.Range("M19").Value = 17.34650563395249
$ws.Range("O19").Value = 24.88726891881454
$ws.Range("B20").Value = 18.40091620006431
$ws.Range("C20").Value = 12.45892557493769
$ws.Range("E20").Value = 15.79836841224538
$ws.Range("F20").Value = 38.42107207422823
$ws.Range("G20").Value = 3.652579920858693
$ws.Range("J20").Value = 8.23459740283838
$ws.Range("L20").Value = 12.36601763841059
$ws.Range("M20").Value = 17.39883288655678
$ws.Range("O20").Value = 24.85330913657645
$ws.Range("B21").Value = 19.01891275112111
$ws.Range("C21").Value = 12.63505764211758
$ws.Range("E21").Value = 15.7557461778397
$ws.Range("F21").Value = 38.35648737953745
$ws.Range("G21").Value = 3.650079805114049
$ws.Range("J21").Value = 8.236476769484238
$ws.Range("L21").Value = 12.39077302930703
$ws.Range("M21").Value = 17.57552914116285
$ws.Range("O21").Value = 24.74711181136979
$ws.Range("B22").Value = 19.41311062621263
$ws.Range("C22").Value = 12.74864252892034
$ws.Range("E22").Value = 15.72931818921997
$ws.Range("F22").Value = 38.32311749785172
$ws.Range("G22").Value = 3.648506473340078
$ws.Range("J22").Value = 8.237961219778136
$ws.Range("L22").Value = 12.4079870465169
$ws.Range("M22").Value = 17.69153841120269
$ws.Range("O22").Value = 24.68362405636963
$ws.Range("B23").Value = 19.20364317153257
$ws.Range("C23").Value = 12.68816699656853
$ws.Range("E23").Value = 15.74329284026957
$ws.Range("F23").Value = 38.34011753167753
$ws.Range("G23").Value = 3.64934065914776
$ws.Range("J23").Value = 8.237145631629078
$ws.Range("L23").Value = 12.39870579048931
$ws.Range("M23").Value = 17.62958215626312
$ws.Range("O23").Value = 24.7169653429292
$ws.Range("B24").Value = 18.39066001338588
$ws.Range("C24").Value = 12.45602260019467
$ws.Range("E24").Value = 15.79908920660837
$ws.Range("F24").Value = 38.42227661426787
$ws.Range("G24").Value = 3.652621811564523
$ws.Range("J24").Value = 8.234571039477039
$ws.Range("L24").Value = 12.36563069064779
$ws.Range("M24").Value = 17.39595566605789
$ws.Range("O24").Value = 24.85514438704014
$ws.Range("B25").Value = 17.47606930264221
$ws.Range("C25").Value = 12.19979246248336
$ws.Range("E25").Value = 15.86542731236978
$ws.Range("F25").Value = 38.5487040909973
$ws.Range("G25").Value = 3.656423190352637
$ws.Range("J25").Value = 8.232906886145031
$ws.Range("L25").Value = 12.33448307362909
$ws.Range("M25").Value = 17.14706397074226
$ws.Range("O25").Value = 25.02941864064746
